$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 15:14:13"
$wsZhCn.Range("H2").Value = "2016-03-22 15:14:34"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 15:14:17"
$wsDeDe.Range("H2").Value = "2016-03-22 15:14:40"
